$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.766.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "'2.674.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'601.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").Value = "'157.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.10%  "

$ws.Range("D9").Value = "'0.124"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "'5.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "'29.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").Value = "'0.0000198"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").Value = "'3.152.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").Value = "'65.557.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "'2.676.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").Value = "'12.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").Value = "'7.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "

$ws.Range("D21").Value = "'352.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.47%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("D24").Value = "'0.0000111"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.80%  "

$ws.Range("E25").Value = "  +4.87%  "

$ws.Range("E26").Value = "  -3.89%  "

$ws.Range("E27").Value = "  +2.14%  "

$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").Value = "'8.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "

$ws.Range("D30").Value = "'544.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.44%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").Value = "'6.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.31%  "

$ws.Range("D35").Value = "'5.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.86%  "

$ws.Range("D36").Value = "'0.423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("D37").Value = "'20.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("E39").Value = "  -0.76%  "

$ws.Range("D40").Value = "'158.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "'42.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "

$ws.Range("D43").Value = "'165.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("D45").Value = "'0.0616"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "

$ws.Range("D46").Value = "'2.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.60%  "

$ws.Range("D47").Value = "'23.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.06%  "

$ws.Range("D48").Value = "'0.646"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("E50").Value = "  +3.50%  "

$ws.Range("D51").Value = "'20.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
